# "maj cours 5 ipt" - mark TP01/TP02 attendance (column C / column D) for
# several students on sheet "S1", move the app window a bit, and update
# the last active selection.

$wb = $excel.ActiveWorkbook

# --- Update workbook window position (bookViews/workbookView in workbook.xml) ---
$win = $excel.ActiveWindow
$win.Left = 6200
$win.Top = 1000

# --- Work on sheet "S1" which holds the attendance/grades table ---
$ws = $wb.Worksheets.Item("S1")
$ws.Activate()

# Mark column D ("TP02") with 1 for the rows that already have a 1 in column C
$ws.Range("D2:D10").Value = 1
$ws.Range("D11:D16").Value = 1

# Rows 11 and 13 were also missing their column C ("TP01") mark
$ws.Range("C11").Value = 1
$ws.Range("C13").Value = 1

# Mark column C ("TP01") with 1 for additional rows (group B)
$ws.Range("C24").Value = 1
$ws.Range("C25").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("C28").Value = 1

# Mark column C ("TP01") with 1 for additional rows (group C)
$ws.Range("C36").Value = 1
$ws.Range("C37").Value = 1
$ws.Range("C38").Value = 1

# Update the active selection on the sheet to E23
$ws.Range("E23").Select()
